$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = "Link polityki prywatności ma poprawny odnośnik"
$ws.Range("E9").Value = "Odnośnik ""polityka prywatności"" zawiera odniesienie do ""polityka-prywatnosci"""

$ws.Range("G9").Select()
